$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Valor Mora" total and the "Cant. Periodos" count ---
$ws.Range("E11").Value = 767563
$ws.Range("F13").Value = 10

# --- Insert two new rows after row 29 (new period 2509 for the two
#     permanent employees), shifting the signature block down ---
$ws.Rows("30:31").Insert(-4121, 0)

# Row 31 (the new last row) gets the "last row" bottom-border style that
# row 29 used to have.
$ws.Range("B29:J29").Copy()
$ws.Range("B31:J31").PasteSpecial(-4122)

# Row 29 becomes a normal "middle" row (copy format from row 28).
$ws.Range("B28:J28").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)

# Row 30 (new) is also a normal "middle" row.
$ws.Range("B28:J28").Copy()
$ws.Range("B30:J30").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# --- Fill in the data for the new period 2509 ---
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "64559974"
$ws.Range("D30").Value = "LIBYS LUZ LOPEZ CAMPILLO"
$ws.Range("E30").Value = "2509"
$ws.Range("F30").Value = 56940
$ws.Range("G30").Value = 1423500

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1043651189"
$ws.Range("D31").Value = "VALENTINA GIRALDO ALVARADO"
$ws.Range("E31").Value = "2509"
$ws.Range("F31").Value = 52000
$ws.Range("G31").Value = 1300000

# --- Center the "Periodo Mora" column for the whole data table ---
$ws.Range("E16:E31").HorizontalAlignment = -4108

Write-Host "done"
